$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing row 55 (2025-03) figures
$ws.Range("B55").Value = 143
$ws.Range("C55").Value = 232
$ws.Range("D55").Value = 61.63793103448276

# Update existing row 56 (2025-04) figures
$ws.Range("B56").Value = 137
$ws.Range("C56").Value = 204
$ws.Range("D56").Value = 67.15686274509804

# Add new row 57 (2025-05)
$ws.Range("A57").Value = "2025-05"
$ws.Range("B57").Value = 36
$ws.Range("C57").Value = 220
$ws.Range("D57").Value = 16.36363636363636
